# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The "CasesTab" row query (cell B2) is updated to drop the optional
# Cohort match/return clause (the blank line right after the first MATCH
# is also removed to line up with the new authored text). Row 2's wrapped
# height shrinks accordingly, and the active selection / scroll position
# moves from B4 up to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and demo.neutered_indicator in [ 'No'] OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Update the Cypher query stored for the "CasesTab" row.
$ws.Range("B2").Value2 = $newCasesQuery

# The text got shorter (one less wrapped line), so the row shrinks from
# 304.5pt to 290pt, matching the height already used by rows 3 and 4.
$ws.Rows(2).RowHeight = 290

# Move the active selection / view up from B4 to B2.
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
